$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the newly-reported tasks for 31/11 (Nhung/Cong/Lan columns, row 10 is a new task row)
$ws.Range("D8").Value = "sửa lỗi đăng nhập ko thành công"
$ws.Range("C9").Value = "sửa lỗi giới tính ở trang đăng nhập"
$ws.Range("D9").Value = "sửa lỗi trang more-info ko hiện thông tin gì"
$ws.Range("B10").Value = "test"
$ws.Range("C10").Value = "làm cái bấm vào điểm đến yêu thích ra list các tour"

# Extend the date merge (A8:A9 -> A8:A10) to cover the new task row, and fix up
# the interior borders so the merged block still reads as one bordered box:
# A9 becomes a middle row (no top/bottom border), A10 becomes the new bottom row
# (top open, bottom border restored).
$ws.Range("A8:A10").Merge()

$a9 = $ws.Range("A9")
$a9.Borders.Item(9).LineStyle = -4142

# Move the active selection like the author left it.
$ws.Range("C14").Select()
